# Updated legacy GSC export data.
# The first daily row (2025-09-03) on the "Chart" sheet is removed; every
# subsequent day's row shifts up by one, so the table now runs one row
# shorter (A1:D81 instead of A1:D82).

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Chart")

# Delete the entire second row (the 2025-09-03 data row); Excel shifts
# everything below it up by one row automatically.
$ws.Rows(2).Delete()
